$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 186 (shifts existing rows 186:307 down to 187:308)
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new record's data.
# (A,B,C,E,F,G,H,I,K,L,M,N,O,P,Q,R mirror the record that used to occupy row 186;
#  D (Fecha) and J (Volumen) are the new values.)
$ws.Range("A186").Value = 3
$ws.Range("B186").Value = "Femacal de La Calera"
$ws.Range("C186").Value = "Coquimbo"
$ws.Range("D186").Value = 44488
$ws.Range("E186").Value = 5
$ws.Range("F186").Value = 100112037
$ws.Range("G186").Value = "Cebollín"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 180
$ws.Range("K186").Value = 3000
$ws.Range("L186").Value = 3000
$ws.Range("M186").Value = 3000
$ws.Range("N186").Value = "$/paquete 36 unidades"
$ws.Range("O186").Value = "Provincia de Quillota"
$ws.Range("P186").Value = 83
$ws.Range("Q186").Value = 36
$ws.Range("R186").Value = "Hortaliza"

# Make sure the date column keeps its original date style/number format.
$ws.Range("D186").NumberFormat = $ws.Range("D187").NumberFormat
